$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New shared strings must be introduced (first referenced) in a precise order
# so the resulting sharedStrings.xml table matches the target file exactly
# (new indices 44-57). We therefore fill the new rows column-by-column in a
# carefully chosen sequence rather than row-by-row or column-by-column.
# ---------------------------------------------------------------------------

# -- introduces shared string 44 "NI 1588-2008 Network Management"
$ws.Range("B28").Value = "NI 1588-2008 Network Management"
# -- introduces shared string 45 "15.0.0"
$ws.Range("C28").Value = "15.0.0"
# -- introduces shared string 46 "NI I/O Trace"
$ws.Range("B29").Value = "NI I/O Trace"
$ws.Range("C29").Value = "15.0"
# -- introduces shared string 47 "2015*"
$ws.Range("E28").Value = "2015*"

# -- introduces shared string 48 "NI PXI Platform Services Configuration"
$ws.Range("B30").Value = "NI PXI Platform Services Configuration"
$ws.Range("C30").Value = "15.0"
$ws.Range("E29").Value = "2015*"
$ws.Range("E30").Value = "2015*"
$ws.Range("B31").Value = "NI PXI Platform Services Configuration"

# -- introduces shared string 49 "2019*"
$ws.Range("E31").Value = "2019*"
# -- introduces shared string 50 "19.5"
$ws.Range("C31").Value = "19.5"

# -- introduces shared string 51 "NI R Series Multifunction RIO"
$ws.Range("B32").Value = "NI R Series Multifunction RIO"
$ws.Range("C32").Value = "15.0"
$ws.Range("E32").Value = "2015*"

# -- introduces shared string 52 "NI Script Editor"
$ws.Range("B33").Value = "NI Script Editor"
# -- introduces shared string 53 "14.0"
$ws.Range("C33").Value = "14.0"
# -- introduces shared string 54 "2014*"
$ws.Range("E33").Value = "2014*"

# -- introduces shared string 55 "NI SignalExpress"
$ws.Range("B34").Value = "NI SignalExpress"
$ws.Range("C34").Value = "2015"
$ws.Range("E34").Value = "2015"

# -- introduces shared string 56 "NI System Configuration"
$ws.Range("B35").Value = "NI System Configuration"
$ws.Range("C35").Value = "19.5"
$ws.Range("E35").Value = "2019*"

# -- introduces shared string 57 "NI-488.2"
$ws.Range("B36").Value = "NI-488.2"
$ws.Range("C36").Value = "15.0"
$ws.Range("E36").Value = "2015*"

$ws.Range("B37").Value = "NI-488.2"
$ws.Range("C37").Value = "19.5"
$ws.Range("E37").Value = "2019*"

# ---------------------------------------------------------------------------
# Fill in the rest (A, D, F columns + G formulas) of the new rows 28-37.
# ---------------------------------------------------------------------------
foreach ($r in 28..37) {
    $ws.Range("A$r").Value = "Corey"
    $ws.Range("D$r").Value = "32"
    $ws.Range("F$r").Value = "Development"
}

# Shared formula across the new rows: =B#&" "&C#
$ws.Range("G28:G37").Formula = "=B28&"" ""&C28"

# ---------------------------------------------------------------------------
# Existing-row edits called out by the diff (reuse shared string 47 "2015*").
# ---------------------------------------------------------------------------

# Row 2: add version-year column + update G2 formula/value
$ws.Range("E2").Value = "2015*"
$ws.Range("G2").Formula = "=B2&"" ""&C2"

# Row 23: add version-year column + update G23 formula/value
$ws.Range("E23").Value = "2015*"
$ws.Range("G23").Formula = "=B23&"" ""&C23"

# Rows 25-27: add the missing "Corey" user cell in column A
$ws.Range("A25").Value = "Corey"
$ws.Range("A26").Value = "Corey"
$ws.Range("A27").Value = "Corey"

# ---------------------------------------------------------------------------
# Column width tweaks. (The COM layer quantizes stored width to whole-pixel
# steps of 1/6 character width, so we pick the input that lands on the
# closest achievable stored width to the target: 5.85546875 -> 5.83333...,
# 33.85546875 -> 33.83333...)
# ---------------------------------------------------------------------------
$ws.Range("E1").EntireColumn.ColumnWidth = 5.0
$ws.Range("G1").EntireColumn.ColumnWidth = 33.0

# ---------------------------------------------------------------------------
# Selection moves to B40 as recorded in the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("B40").Select()
